# games_one.xlsx edit ("Add files via upload"):
#  - games!G4/H4 get the same lat/lng text values already used for the first
#    "BP One Teszt" task row (tasks!C12:D12, shared strings 64/65), and
#    games!I4 gets the matching geofenceRadius of 10 (mirrors rows 2 & 3).
#  - The "games" sheet becomes the active/selected sheet with I4 selected.
#  - The "tasks" sheet selection moves to C12:D12 and it stops being the
#    active tab.

$wb = $excel.ActiveWorkbook
$wsGames = $wb.Worksheets.Item("games")
$wsTasks = $wb.Worksheets.Item("tasks")

# --- Fill in the missing location/radius cells on row 4 of "games" ---
# Copy the existing text cells from tasks!C12:D12 so the new games!G4/H4
# cells reuse the very same shared-string values (rather than becoming new
# numeric cells), exactly like the source edit.
$wsTasks.Range("C12:D12").Copy() | Out-Null
$wsGames.Range("G4").PasteSpecial(-4163) | Out-Null
$wsGames.Range("I4").Value = 10

# --- Update selections / active sheet to match the saved view state ---
$wsTasks.Activate() | Out-Null
$wsTasks.Range("C12:D12").Select() | Out-Null

$wsGames.Activate() | Out-Null
$wsGames.Range("I4").Select() | Out-Null
